$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New data row (row 2): identifier, title, levelOfDescription, extentAndMedium, notes
# Match the font used for the rest of the new row's cells (Calibri 10pt, theme color 1)
$cells = @("A2","C2","D2","E2","F2","G2","H2")
foreach ($addr in $cells) {
    $c = $ws.Range($addr)
    $c.Font.Name = "Calibri"
    $c.Font.Size = 10
    $c.Font.ThemeColor = 1
}

$ws.Range("A2").Value = "MCH142"
$ws.Range("C2").Value = "OBSTACLES TO NEGOTIATION"
$ws.Range("E2").Value = "Series"
$ws.Range("F2").Value = "1 Box"
$ws.Range("G2").Value = "LOCATION: CABINET 1B | GRAP COUNT NUMER: NONE"

# Restore/refresh the frozen header pane with the new row selected.
$ws.Range("A2").Select() | Out-Null
$excel.ActiveWindow.FreezePanes = $true
$ws.Range("A2:I2").Select() | Out-Null
